$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "weight" column header in P1 (becomes a new shared string), matching the
# style already used by the neighbouring header cells in row 1.
$ws.Cells.Item(1, 16).Value = "weight"

# Fill P2:P30 with the constant weight value (20) used for every data row.
$ws.Range("P2:P30").Value = 20

# Match the style used by the adjacent column O cells (same default numeric style).
$ws.Range("P1:P30").Style = $ws.Range("O1:O30").Style

# Scroll the view so column B is the left-most visible column, then leave the
# new P2:P30 block selected - this is the view state the author ended up with
# after filling in the new column.
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("P2:P30").Select()
